# Automatic update of files.
#
# The underlying data rows (2-17) get re-shuffled: each destination row's
# full contents (every column A:AY) become equal to some *other* row's
# original contents, per a fixed permutation recovered from the diff.
# Row 1 (headers) is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# destination row -> source row (source row's ORIGINAL content moves into
# destination row)
$perm = @{
    2  = 12
    3  = 11
    4  = 14
    5  = 15
    6  = 10
    7  = 17
    8  = 9
    9  = 6
    10 = 3
    11 = 4
    12 = 8
    13 = 7
    14 = 16
    15 = 13
    16 = 2
    17 = 5
}

# Snapshot every source row (2-17, columns A:AY) BEFORE any writes, since the
# permutation is not its own inverse and rows would otherwise clobber each
# other's source data mid-script. Also remember the literal text of the
# "date-looking" columns (Y = Startdatum, AA = Slutdatum) separately so we
# can re-assert them as text afterwards (a plain Value2 array write lets
# Excel auto-coerce "2023-09-03" into a real date serial, which we don't
# want - the source file stores it as literal text).
$snapshot = @{}
$snapshotY = @{}
$snapshotAA = @{}
for ($r = 2; $r -le 17; $r++) {
    $snapshot[$r] = $ws.Range("A" + $r + ":AY" + $r).Value2
    $snapshotY[$r] = $ws.Range("Y" + $r).Value2
    $snapshotAA[$r] = $ws.Range("AA" + $r).Value2
}

# Now write each destination row from the snapshot of its mapped source row.
foreach ($destRow in $perm.Keys) {
    $srcRow = $perm[$destRow]
    $ws.Range("A" + $destRow + ":AY" + $destRow).Value2 = $snapshot[$srcRow]
    # Re-assert Y/AA as literal text (apostrophe prefix forces text without
    # altering the cell's number format/style).
    $ws.Range("Y" + $destRow).Value = "'" + $snapshotY[$srcRow]
    $ws.Range("AA" + $destRow).Value = "'" + $snapshotAA[$srcRow]
}
